$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F and G on row 1
$ws.Range("F1").Value = "Code Desc"
$ws.Range("G1").Value = "Level5"

# New coding hierarchy values for row 213 (G, H, I)
$ws.Range("G213").Value = "Rental"
$ws.Range("H213").Value = "Transport"
$ws.Range("I213").Value = "Leisure"

# Update the active selection to match the author's final cursor position
$ws.Range("G227").Select()
